$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jake McGee's injury report (row 5) has been removed from the newsletter.
$ws.Rows.Item(5).Delete() | Out-Null

# Kyle Freeland's entry (now row 4) gets an updated "Last.Updated" date and
# a refreshed injury-details blurb.
$ws.Range("C4").Value = "August 09 2017"
$ws.Range("E4").Value = "Freeland is on the 10-day disabled list with a left groin strain and is expected to remain sidelined until the end of August."

# The custom/bold font formatting that used to highlight the "Name" column
# (A2:A4) has been removed - revert those cells back to the default style.
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Style = "Normal"

# Reflect the cursor position saved in the workbook after the edit.
$ws.Range("A18").Select() | Out-Null
